$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Swap the "C. Kata" / "Cia M." rows (originally rows 3-4) with the
# "Charles Y." / "Nate M." rows (originally rows 5-6). The two Perkinite
# entries (C. Kata / Cia M.) move down to rows 5-6 and get re-classed (CS3)
# which is reflected by them picking up fresh direct formatting rather than
# the alternating table-band fill used by the rest of the rows.
# ---------------------------------------------------------------------------

# Row 3 (was "C. Kata") becomes "Charles Y."
$ws.Range("A3").Value = "Charles Y."
$ws.Range("B3").Value = 75
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = "Shinai"

# Row 4 (was "Cia M.") becomes "Nate M."
$ws.Range("A4").Value = "Nate M."
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = "Claws"

# Row 5 (was "Charles Y.") becomes "C. Kata"
$ws.Range("A5").Value = "C. Kata"
$ws.Range("B5").Value = 75
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 17
$ws.Range("F5").Value = "Railgun"

# Row 6 (was "Nate M.") becomes "Cia M."
$ws.Range("A6").Value = "Cia M."
$ws.Range("B6").Value = 60
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 17
$ws.Range("F6").Value = "Magic Wand"

# ---------------------------------------------------------------------------
# Re-format rows 5 and 6: both keep the light-blue row fill, but they no
# longer follow the regular odd/even band - row 5 keeps the "interior" thin
# white borders, row 6 picks up the "last row" border treatment (no bottom
# border), matching the look applied when the table was touched up.
# ---------------------------------------------------------------------------

$lightBlue = 15853276
$white = 16777215

# Row 5: interior-row border treatment (thin white right+bottom on A:E, thin white bottom on F)
$r5 = $ws.Range("A5:F5")
$r5.Interior.Color = $lightBlue

$r5left = $ws.Range("A5:E5")
$r5left.Borders.Item(10).LineStyle = 1
$r5left.Borders.Item(10).Weight = 2
$r5left.Borders.Item(10).Color = $white

$r5.Borders.Item(9).LineStyle = 1
$r5.Borders.Item(9).Weight = 2
$r5.Borders.Item(9).Color = $white

# Row 6: last-row border treatment (thin white right only on A:E, no bottom border)
$r6 = $ws.Range("A6:F6")
$r6.Interior.Color = $lightBlue

$r6left = $ws.Range("A6:E6")
$r6left.Borders.Item(10).LineStyle = 1
$r6left.Borders.Item(10).Weight = 2
$r6left.Borders.Item(10).Color = $white

$r6.Borders.Item(9).LineStyle = -4142

# ---------------------------------------------------------------------------
# Column A now gets an explicit width, and the remembered selection moves.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 11.6

$ws.Range("H8").Select()
